$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5, 4, "2025-04-26 11:22:52", "John Smith found battery 2. Now John Smith is Tired"),
    @(6, 5, "2025-04-26 11:24:48", "John Smith found battery 3. Now John Smith is Confident"),
    @(7, 6, "2025-04-26 11:32:09", "John Smith found battery 2. Now John Smith is Confident"),
    @(8, 7, "2025-04-26 11:40:24", "John Smith found battery 2. Now John Smith is Happy`n"),
    @(9, 8, "2025-04-26 11:55:49", "John Smith found battery 1. `n Now John Smith is Frustrated`n"),
    @(10, 9, "2025-04-26 11:59:03", "John Smith found battery 2.`nNow John Smith is Tired`n")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 3).WrapText = $true
}
